$d = $word.ActiveDocument

# Replace the semicolon with a colon after "Samtykkeerklæring"
$d.Content.Find.Execute("Samtykkeerklæring;", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Samtykkeerklæring:", 2)
